# Improved accuracy of stimulus presentation time-logging
# Renames task-order sheets with updated timestamps and refreshes the
# per-trial stimulus file names (and eyes-open/closed ordering) with new
# timestamped values.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555109597478"
$ws1.Range("B2").Value = "go_stims-16512555109257472.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555109437459.csv"
$ws1.Range("B4").Value = "go_stims-16512555109447474.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555109587457.csv"

# --- Sheet 2: NB_TO ------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16512555130667443"
$ws2.Range("B2").Value = "ZB-match_2-16512555117517457.csv"
$ws2.Range("B3").Value = "OB-16512555120587444.csv"
$ws2.Range("B4").Value = "OB-16512555121307464.csv"
$ws2.Range("B5").Value = "TB-16512555125597453.csv"
$ws2.Range("B6").Value = "ZB-match_5-16512555118897471.csv"
$ws2.Range("B7").Value = "OB-16512555119097455.csv"
$ws2.Range("B8").Value = "TB-16512555130507455.csv"
$ws2.Range("B9").Value = "TB-16512555130047457.csv"
$ws2.Range("B10").Value = "ZB-match_6-1651255511185746.csv"

# --- Sheet 3: RS_TO -------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651255513073745"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO -------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512555131317508"
$ws4.Range("B2").Value = "MM_stims-1651255513097748.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255513075749.csv"
$ws4.Range("B4").Value = "MM_stims-16512555131137452.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555130987487.csv"
$ws4.Range("B6").Value = "MM_stims-1651255513129746.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555131147473.csv"

# --- Sheet 5: vSAT_TO -------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512555132098255"
$ws5.Range("B2").Value = "SAT_stims-16512555131617472.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651255513193748.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555131367471.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555131777444.csv"
